$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

$ws.Range("B2").Value = 2539.55
$ws.Range("C2").Value = 2530.3
$ws.Range("B3").Value = 395.6
$ws.Range("C3").Value = 393.35
$ws.Range("B4").Value = 1468.6
$ws.Range("C4").Value = 1470.1
$ws.Range("B5").Value = 7137.45
$ws.Range("C5").Value = 7116.45
$ws.Range("B6").Value = 234.9
$ws.Range("C6").Value = 237.25
$ws.Range("B7").Value = 191.85
$ws.Range("C7").Value = 193.35
$ws.Range("B8").Value = 44559.1
$ws.Range("C8").Value = 44519.45
$ws.Range("B9").Value = 483.55
$ws.Range("C9").Value = 480.7
$ws.Range("B10").Value = 3341.65
$ws.Range("C10").Value = 3358
$ws.Range("B11").Value = 140.7
$ws.Range("C11").Value = 141.7
$ws.Range("B12").Value = 1171.25
$ws.Range("C12").Value = 1182.5
$ws.Range("B13").Value = 1425.15
$ws.Range("C13").Value = 1401
$ws.Range("B14").Value = 649.1
$ws.Range("C14").Value = 660.75
$ws.Range("B15").Value = 421.55
$ws.Range("C15").Value = 423
$ws.Range("B16").Value = 1542.2
$ws.Range("C16").Value = 1542.45
$ws.Range("B17").Value = 304.05
$ws.Range("C17").Value = 298.05
$ws.Range("B18").Value = 19395.6
$ws.Range("C18").Value = 19430.9
$ws.Range("B19").Value = 576.95
$ws.Range("C19").Value = 576.9
$ws.Range("B20").Value = 582.55
$ws.Range("C20").Value = 568.25
$ws.Range("B21").Value = 610.95
$ws.Range("C21").Value = 613.9
$ws.Range("B22").Value = 249.75
$ws.Range("C22").Value = 245.6
$ws.Range("B23").Value = 118.25
$ws.Range("C23").Value = 118.95
